# Update the "取得日時" (retrieved datetime) timestamp stored in column A
# of the "ランサーズ" sheet for every data row (rows 2-13), changing the
# old timestamp "2025-12-25 12:37:44" to the new one "2025-12-25 12:50:27".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-12-25 12:37:44"
$newTimestamp = "2025-12-25 12:50:27"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 13
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
